$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix up the timestamp recorded for the previous "Actualizar" run
#    (rows 380-393): 44232.24014104684 -> 44232.24014104166
# ---------------------------------------------------------------------------
for ($r = 380; $r -le 393; $r++) {
    $ws.Cells.Item($r, 4).Value2 = 44232.24014104166
}

# ---------------------------------------------------------------------------
# 2) Append the new "Actualizar" run: rows 394-407
# ---------------------------------------------------------------------------
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$urls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
# the MapStore row (index 8, row 402) links with a trailing "/" fragment;
# its visible cell text already carries that fragment (matches the text
# used by every earlier MapStore row in the sheet).
$subAddresses = @($null,$null,$null,$null,$null,$null,$null,$null,"/",$null,$null,$null,$null,$null)
$displayText  = @()
for ($i = 0; $i -lt $urls.Length; $i++) {
    if ($subAddresses[$i]) {
        $displayText += ($urls[$i] + "#" + $subAddresses[$i])
    } else {
        $displayText += $urls[$i]
    }
}

$startRow = 394
$newDate = 44232.26129783387

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value = $names[$i]

    # Set the display text first so Hyperlinks.Add's own text reset keeps
    # re-using the existing shared string instead of minting a new one.
    $ws.Cells.Item($r, 2).Value = $displayText[$i]
    if ($subAddresses[$i]) {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $urls[$i], $subAddresses[$i]) | Out-Null
    } else {
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $urls[$i]) | Out-Null
    }
    # Hyperlinks.Add mints its own (duplicate) cell format; snap the cell
    # back onto the sheet's single shared "Hyperlink" style.
    $ws.Cells.Item($r, 2).Style = "Hyperlink"

    $ws.Cells.Item($r, 3).Value = "Disponible"

    $ws.Cells.Item($r, 4).Value2 = $newDate
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
